$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Bitcoin USD (row 4) close price and 5-day return
$ws.Range("D4").Value = 90865.73
$ws.Range("F4").Value = 4.68

# Update MACRO_SCORE column (N) for all data rows
$ws.Range("N2").Value = 85.83574689470727
$ws.Range("N3").Value = 85.83574689470727
$ws.Range("N4").Value = 85.83574689470727
$ws.Range("N5").Value = 85.83574689470727
$ws.Range("N6").Value = 85.83574689470727
